# Apply "Add data for 2022-06-06" update to the carjacking-by-month-yoy workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-05-29"

# 2. Update the label text for the May row (A6) to reflect the new date.
$ws.Range("A6").Value = "May (through 05-29)"

# 3. Update the May row values (row 6).
$ws.Range("B6").Value = 18
$ws.Range("D6").Value = 56
$ws.Range("G6").Value = 64
$ws.Range("H6").Value = 103
$ws.Range("I6").Value = 107

# 4. Update the Total row values (row 7) to match the new sums.
$ws.Range("B7").Value = 107
$ws.Range("D7").Value = 309
$ws.Range("G7").Value = 326
$ws.Range("H7").Value = 626
$ws.Range("I7").Value = 658
